$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 1.459612070389937;  C = 114.8270160096505;  D = 26.21740644021617;  E = 8.660232485948974;   G = 151.1642670062056 }
    3  = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 0.8054896365839992; E = 0.496779210170732;   G = 6.201049113329182 }
    4  = @{ B = 0.127881588408715;  C = 1.667794583268128;   D = 0.8054896365839992; E = 645.3272768299601;   G = 647.9284426382209 }
    5  = @{ B = 1.459612070389937;  C = 1.667794583268128;   D = 3.900430680208489;  E = 0.496779210170732;   G = 7.524616544037286 }
    6  = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 0.1575252929769615; E = 0.496779210170732;   G = 5.553084769722144 }
    7  = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 26.21740644021617;  E = 0.496779210170732;   G = 31.61296591696135 }
    8  = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 0.8054896365839992; E = 0.496779210170732;   G = 6.201049113329182 }
    9  = @{ B = 1.459612070389937;  C = 0.3127903958511391;  D = 0.8054896365839992; E = 0.496779210170732;   G = 3.074671312995807 }
    10 = @{ B = 3.230985683306322;  C = 1.667794583268128;   D = 3.900430680208489;  E = 0.496779210170732;   G = 9.295990156953671 }
    11 = @{ B = 1.459612070389937;  C = 1.667794583268128;   D = 3.900430680208489;  E = 0.496779210170732;   G = 7.524616544037286 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("C$row").Value = $vals.C
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
    $ws.Range("G$row").Value = $vals.G
}
